$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record above the current row 159, shifting the
# existing rows (159-165) down to (160-166).
$ws.Rows.Item(159).Insert()

# Populate the new row 159 with this week's Aji price entry.
$ws.Range("A159").Value = 4
$ws.Range("B159").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C159").Value = "Los Lagos"
$ws.Range("D159").Value = 44509
$ws.Range("E159").Value = 10
$ws.Range("F159").Value = 100112021
$ws.Range("G159").Value = "Ají"
$ws.Range("H159").Value = "Inferno"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 80
$ws.Range("K159").Value = 35000
$ws.Range("L159").Value = 35000
$ws.Range("M159").Value = 35000
$ws.Range("N159").Value = "$/caja 12 kilos"
$ws.Range("O159").Value = "Región de Arica y Parinacota"
$ws.Range("P159").Value = 2917
$ws.Range("Q159").Value = 12
$ws.Range("R159").Value = "Hortaliza"
